$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column (H1), reusing the header style already
# used by the other column headers (e.g. G1) via a format copy.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the new data values (plain numeric cells, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
